# Update the transition-probability matrix on Sheet1 with newly recomputed
# values (more simulated games -> updated empirical probabilities).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.3095238095238095
$ws.Cells.Item(2, 3).Value = 0.3095238095238095
$ws.Cells.Item(2, 16).Value = 0.3095238095238095
$ws.Cells.Item(2, 19).Value = 0.07142857142857142

$ws.Cells.Item(3, 16).Value = 0.9230769230769231
$ws.Cells.Item(3, 19).Value = 0.07692307692307693

$ws.Cells.Item(4, 10).Value = 0.3333333333333333
$ws.Cells.Item(4, 16).Value = 0.6666666666666666

$ws.Cells.Item(6, 6).Value = 0.1538461538461539
$ws.Cells.Item(6, 10).Value = 0.3076923076923077
$ws.Cells.Item(6, 17).Value = 0.3076923076923077
$ws.Cells.Item(6, 18).Value = 0.07692307692307693
$ws.Cells.Item(6, 19).Value = 0.1538461538461539

$ws.Cells.Item(7, 2).Value = 0.09090909090909091
$ws.Cells.Item(7, 10).Value = 0.3636363636363636
$ws.Cells.Item(7, 17).Value = 0.2727272727272727
$ws.Cells.Item(7, 18).Value = 0.09090909090909091
$ws.Cells.Item(7, 19).Value = 0.1818181818181818

$ws.Cells.Item(8, 2).Value = 0.2222222222222222
$ws.Cells.Item(8, 10).Value = 0.1111111111111111
$ws.Cells.Item(8, 17).Value = 0.3703703703703703
$ws.Cells.Item(8, 18).Value = 0.03703703703703703
$ws.Cells.Item(8, 19).Value = 0.2592592592592592

$ws.Cells.Item(9, 6).Value = 0.1111111111111111
$ws.Cells.Item(9, 10).Value = 0.2222222222222222
$ws.Cells.Item(9, 17).Value = 0.1111111111111111
$ws.Cells.Item(9, 19).Value = 0.5555555555555556

$ws.Cells.Item(10, 2).Value = 0.1833333333333333
$ws.Cells.Item(10, 4).Value = 0.025
$ws.Cells.Item(10, 6).Value = 0.04166666666666666
$ws.Cells.Item(10, 10).Value = 0.1833333333333333
$ws.Cells.Item(10, 17).Value = 0.1666666666666667
$ws.Cells.Item(10, 18).Value = 0.075
$ws.Cells.Item(10, 19).Value = 0.325

$ws.Cells.Item(11, 7).Value = 0.1764705882352941
$ws.Cells.Item(11, 11).Value = 0.1764705882352941
$ws.Cells.Item(11, 12).Value = 0.6470588235294118

$ws.Cells.Item(12, 7).Value = 0.5454545454545454
$ws.Cells.Item(12, 10).Value = 0.4545454545454545

$ws.Cells.Item(15, 8).Value = 0.2
$ws.Cells.Item(15, 9).Value = 0.06666666666666667
$ws.Cells.Item(15, 10).Value = 0.4666666666666667
$ws.Cells.Item(15, 11).Value = 0.1333333333333333
$ws.Cells.Item(15, 15).Value = 0.06666666666666667
$ws.Cells.Item(15, 19).Value = 0.06666666666666667

$ws.Cells.Item(16, 6).Value = 0.03846153846153846
$ws.Cells.Item(16, 8).Value = 0.1538461538461539
$ws.Cells.Item(16, 9).Value = 0.07692307692307693
$ws.Cells.Item(16, 10).Value = 0.5384615384615384
$ws.Cells.Item(16, 11).Value = 0.03846153846153846
$ws.Cells.Item(16, 15).Value = 0.03846153846153846
$ws.Cells.Item(16, 19).Value = 0.1153846153846154

$ws.Cells.Item(17, 6).Value = 0.02631578947368421
$ws.Cells.Item(17, 8).Value = 0.1052631578947368
$ws.Cells.Item(17, 9).Value = 0.05263157894736842
$ws.Cells.Item(17, 10).Value = 0.4736842105263158
$ws.Cells.Item(17, 11).Value = 0.1052631578947368
$ws.Cells.Item(17, 13).Value = 0.02631578947368421
$ws.Cells.Item(17, 15).Value = 0.1052631578947368
$ws.Cells.Item(17, 19).Value = 0.1052631578947368

$ws.Cells.Item(18, 6).Value = 0.08333333333333333
$ws.Cells.Item(18, 8).Value = 0.25
$ws.Cells.Item(18, 9).Value = 0.08333333333333333
$ws.Cells.Item(18, 10).Value = 0.3333333333333333
$ws.Cells.Item(18, 13).Value = 0.08333333333333333
$ws.Cells.Item(18, 15).Value = 0.08333333333333333
$ws.Cells.Item(18, 19).Value = 0.08333333333333333

$ws.Cells.Item(19, 9).Value = 0.05128205128205128
$ws.Cells.Item(19, 10).Value = 0.4615384615384616
$ws.Cells.Item(19, 11).Value = 0.08974358974358974
$ws.Cells.Item(19, 15).Value = 0.07692307692307693
$ws.Cells.Item(19, 19).Value = 0.1538461538461539
